$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Last Updated: 26 February 2026" -> "Last Updated: 27 February 2026"
#    The target OOXML keeps this split across three runs:
#      "Last Updated: 2" | "7" | " February 2026"
#    so edit it character-by-character rather than via a single
#    whole-phrase replace (which would collapse back into one run).
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Last Updated: 26 February 2026", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$luStart = $findRng.Start

# Force a run split around the "6" character (offset 15..16 within the phrase)
# by toggling a character property on and back off - this materializes an
# explicit (empty) <w:rPr/> on the newly separated runs without altering
# their visible formatting.
$digit = $d.Range($luStart + 15, $luStart + 16)
$digit.Bold = 1
$digit.Bold = 0

# Now replace the (still isolated) run's text "6" -> "7".
$digit2 = $d.Range($luStart + 15, $luStart + 16)
$digit2.Text = "7"

# Re-materialize <w:rPr/> on the left/right runs (the text assignment above
# can re-merge formatting-equal neighbours and drop their explicit rPr).
$left = $d.Range($luStart, $luStart + 15)
$left.Bold = 1
$left.Bold = 0
$right = $d.Range($luStart + 16, $luStart + 31)
$right.Bold = 1
$right.Bold = 0

# ------------------------------------------------------------------
# 2) Collapse the "references" sentence (previously split over several
#    runs) into a single run with the corrected wording.
# ------------------------------------------------------------------
$oldRefs = "You will need to obtain the CP/M MSX loader program from references [5] and [10], including the BIOS file msx-us.rom and some MSX games from reference [4] and transfer them to your H-89."
$d.Content.Find.Execute($oldRefs, $true, $false, $false, $false, $false, $true, 1, $false, $oldRefs, 2) | Out-Null

# Materialize an explicit <w:rPr/> on the now-merged run.
$refsRng = $d.Content
$refsRng.Find.Execute($oldRefs, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$refsRng.Bold = 1
$refsRng.Bold = 0

# ------------------------------------------------------------------
# 3) "The base i/o ports (in hex) used..." -> "...(in hex an d octal)..."
# ------------------------------------------------------------------
$d.Content.Find.Execute("The base i/o ports (in hex) used by the board are shown below:", $true, $false, $false, $false, $false, $true, 1, $false, "The base i/o ports (in hex an d octal) used by the board are shown below:", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Port table cells: add the octal equivalent after each hex value.
#    Each cell is (and must remain) a single run carrying an explicit
#    (empty) <w:rPr/>, so re-materialize it after the text swap.
# ------------------------------------------------------------------
function Update-PortCell([string]$oldText, [string]$newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
    $cellRng = $d.Content
    $cellRng.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $cellRng.Bold = 1
    $cellRng.Bold = 0
}

Update-PortCell "`$A0" "`$A0  240Q"
Update-PortCell "`$BA" "`$BA  272Q"
Update-PortCell "`$98" "`$98  230Q"
Update-PortCell "`$B8" "`$B8  270Q"
